$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells we touch stay as text (avoid Excel auto-converting numeric-looking
# strings like "0.999" or "5.50" into actual numbers).

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.135.40'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.02%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.390.12'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -5.13%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.26%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '478.76'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.45%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.63'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.11%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.27%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.499'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.14%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.381.55'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -6.13%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0974'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.44%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.50'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.20%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.325'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.08%  '

# Row 13
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.60%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.798.05'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -5.10%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '56.392.55'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.50%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.38'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.56%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.81%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.380.67'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -5.85%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.46'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.24%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '313.59'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.75%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.75'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.20%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.20%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.67'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.95%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '56.89'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.95%  '

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.41%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.396'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.27%  '

# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.53%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.524.75'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.51%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.29'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.13%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0774'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.64%  '

# Row 31
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.13%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.24'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.77%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.95'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.01%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.48'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.34%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.98'

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.11'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.27%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.856'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.03%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.61'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.87%  '

# Row 39
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.52'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.38%  '

# Row 40
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.35'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.68%  '

# Row 41
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.995'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.22%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0544'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.70%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.38'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.72%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.584'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.03%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0944'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.54%  '

# Row 46
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.19%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '256.75'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.95%  '

# Row 48
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0223'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.75%  '

# Row 49
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.56'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -6.10%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.04'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.85%  '

# Row 51
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.67'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +16.96%  '
